$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheets: 1 = LeetCode, 2 = AlgoExpert, 3 = Blind 150
# ---------------------------------------------------------------------------
$wsLeet  = $wb.Worksheets.Item(1)
$wsAlgo  = $wb.Worksheets.Item(2)
$wsBlind = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# Blind 150 sheet: row 14's Type changes from "Binary Search" to "Trees",
# and two new rows (15, 16) are appended for "Balanced Binary Tree" and
# "Trie". Editing this sheet first so the brand-new shared strings line up
# with the index order seen in the target workbook (121-124 come from this
# sheet, 125-126 come from the AlgoExpert sheet edited afterwards).
# ---------------------------------------------------------------------------
$wsBlind.Activate()

$wsBlind.Range("B14").Value = "Trees"

$wsBlind.Range("A15").Value = 14
$wsBlind.Range("B15").Value = "Trees"
$wsBlind.Range("C15").Value = "Balanced Binary Tree"
$wsBlind.Range("D15").Value = "1. same as maxDepth, just use l - r < 1 before return"
$wsBlind.Rows.Item(15).RowHeight = 28

$wsBlind.Range("A16").Value = 15
$wsBlind.Range("B16").Value = "Trie"
$wsBlind.Range("C16").Value = "Implement a Trie"
$wsBlind.Rows.Item(16).RowHeight = 14

$wsBlind.Range("D16").Select() | Out-Null

# ---------------------------------------------------------------------------
# AlgoExpert sheet: fill in row 7 with the "Longest Peak" problem.
# ---------------------------------------------------------------------------
$wsAlgo.Activate()

$wsAlgo.Range("A7").Value = 6
$wsAlgo.Range("B7").Value = "Medium"
$wsAlgo.Range("C7").Value = "Longest Peak"
$notes = "1. maintain a maxLength, and ptr" + [char]10 + `
  "2. increament inc and dec within the while for i < size both within while of A[i-1]< A[i] and opp" + [char]10 + `
  "3. math,max between maxlength and inc+dec+1" + [char]10 + `
  "4. if pre and curr are equal, just i++" + [char]10 + `
  "5. Pattern of while inc while dec if equal till inc again"
$wsAlgo.Range("D7").Value = $notes
$wsAlgo.Rows.Item(7).RowHeight = 70

$wsAlgo.Range("D17:D22").Select() | Out-Null

# AlgoExpert ends up the active / visible tab.
$wsAlgo.Activate()
